$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace sample rows (4-5: 할아버지/잼민이) with monster rows (4-7: mon1..mon4) ---

# Column A (Key)
$ws.Range("A4").Value = "mon1"
$ws.Range("A5").Value = "mon2"

# Column C (Monster type)
$ws.Range("C4").Value = "Mob"
$ws.Range("C5").Value = "Mob"
$ws.Range("C6").Value = "Mob"
$ws.Range("C7").Value = "Mob"

# Column B (Name / description)
$ws.Range("B4").Value = "mon1_desc"
$ws.Range("B5").Value = "mon2_desc"

# New rows 6-7, column A
$ws.Range("A6").Value = "mon3"
$ws.Range("A7").Value = "mon4"

# New rows 6-7, column B
$ws.Range("B6").Value = "mon3_desc"
$ws.Range("B7").Value = "mon4_desc"

# Clear the old Icon(I) / SpecialMonsterTag(K) sample values that no longer apply
$ws.Range("I4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("I5").ClearContents()

# Column H (Prefabs)
$ws.Range("H4").Value = "mon1"
$ws.Range("H5").Value = "mon2"
$ws.Range("H6").Value = "mon3"
$ws.Range("H7").Value = "mon4"

# Column J (IsSpecial) - all false for the new monster rows
$ws.Range("J4").Value = $false
$ws.Range("J5").Value = $false
$ws.Range("J6").Value = $false
$ws.Range("J7").Value = $false

# Update selection to match the saved state
$ws.Range("A5").Select()
